# Applies a 3-way cyclic rotation of the weekly price data blocks:
#   rows 2-3 (old date 44559) <- values from rows 8-9 (old date 44216)
#   rows 6-7 (old date 44574) <- values from rows 2-3 (old date 44559)
#   rows 8-9 (old date 44216) <- values from rows 6-7 (old date 44574)
# Only columns D (Fecha) and N, O, P, S (prices) change; Q/R/etc. stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original ("before") values for the rows that will move, since
# writes to one block must not affect the source values still needed by
# another block.
function Get-RowData($r) {
    return @{
        D = $ws.Cells.Item($r, 4).Value2    # Fecha
        N = $ws.Cells.Item($r, 14).Value2   # Precio minimo
        O = $ws.Cells.Item($r, 15).Value2   # Precio maximo
        P = $ws.Cells.Item($r, 16).Value2   # Precio promedio ponderado
        S = $ws.Cells.Item($r, 19).Value2   # Precio $/Kg
    }
}

$row2 = Get-RowData 2
$row3 = Get-RowData 3
$row6 = Get-RowData 6
$row7 = Get-RowData 7
$row8 = Get-RowData 8
$row9 = Get-RowData 9

function Set-RowData($r, $data) {
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 14).Value = $data.N
    $ws.Cells.Item($r, 15).Value = $data.O
    $ws.Cells.Item($r, 16).Value = $data.P
    $ws.Cells.Item($r, 19).Value = $data.S
}

# rows 2,3 take what rows 8,9 had
Set-RowData 2 $row8
Set-RowData 3 $row9

# rows 6,7 take what rows 2,3 had
Set-RowData 6 $row2
Set-RowData 7 $row3

# rows 8,9 take what rows 6,7 had
Set-RowData 8 $row6
Set-RowData 9 $row7

$wb.Save()
